# Reorder the "Requisitos" bullet list so the LOB1012 line moves from
# the first entry to the last entry in the list, i.e.:
#   LOB1012 / LOQ4095 / LOQ4098  ->  LOQ4095 / LOQ4098 / LOB1012

$d = $word.ActiveDocument

$lobText = "LOB1012 -  Estatística  (Requisito fraco)"

# Locate the LOB1012 run and extend the range to also capture its
# trailing line break (w:br), which Word represents as a vertical-tab
# character (chr 11) in Range.Text.
$lobRng = $d.Content
$found = $lobRng.Find.Execute($lobText, $true, $false, $false, $false, $false, `
                               $true, 1, $false, "", 0)

if ($found) {
    [void]$lobRng.MoveEnd(1, 1)

    # Remove the "LOB1012 ... (Requisito fraco)" line (text + line break)
    # from its current (first) position in the list.
    [void]$lobRng.Delete()

    # The list is the last paragraph under the "Requisitos" heading; find
    # it again now that the text has shifted, and insert the removed line
    # right before the paragraph mark (i.e. after the last remaining line
    # break), restoring its own trailing line break.
    $paraRng = $d.Paragraphs.Last.Range
    $insertPos = $paraRng.End - 1
    $insertRng = $d.Range($insertPos, $insertPos)
    [void]$insertRng.InsertAfter($lobText + [char]11)
}
